$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.008.49"
$ws.Range("E2").Value = "  +4.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.416.08"
$ws.Range("E3").Value = "  +3.23%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.14"
$ws.Range("E5").Value = "  +3.39%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.66"
$ws.Range("E6").Value = "  +8.24%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.416.69"
$ws.Range("E8").Value = "  +3.07%  "
$ws.Range("E9").Value = "  +2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.50"
$ws.Range("E10").Value = "  +1.41%  "
$ws.Range("E11").Value = "  +9.96%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.395"
$ws.Range("E12").Value = "  +6.76%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.001.11"
$ws.Range("E13").Value = "  +3.99%  "
$ws.Range("E14").Value = "  +2.06%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000180"
$ws.Range("E15").Value = "  +8.85%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.419.72"
$ws.Range("E16").Value = "  +4.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.54"
$ws.Range("E17").Value = "  +6.64%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.042.76"
$ws.Range("E18").Value = "  +4.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.19"
$ws.Range("E19").Value = "  +7.26%  "
$ws.Range("E20").Value = "  +4.85%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.55"
$ws.Range("E21").Value = "  +7.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "391.72"
$ws.Range("E22").Value = "  +11.81%  "
$ws.Range("E23").Value = "  +3.85%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.553.38"
$ws.Range("E24").Value = "  +3.53%  "
$ws.Range("E25").Value = "  +19.19%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "71.53"
$ws.Range("E27").Value = "  +4.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.59"
$ws.Range("E28").Value = "  +11.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.67"
$ws.Range("E29").Value = "  +5.77%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.32"
$ws.Range("E31").Value = "  +6.94%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.161"
$ws.Range("E32").Value = "  +6.92%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.18"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.447.60"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("E35").Value = "  +0.05%  "
$ws.Range("E36").Value = "  +3.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.50"
$ws.Range("E37").Value = "  +4.46%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "7.02"
$ws.Range("E38").Value = "  +3.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  +7.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.57"
$ws.Range("E40").Value = "  +3.47%  "
$ws.Range("E41").Value = "  +6.30%  "
$ws.Range("E42").Value = "  +15.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.793"
$ws.Range("E43").Value = "  +6.58%  "
$ws.Range("B44").Value = "ONDO"
$ws.Range("C44").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.24"
$ws.Range("E44").Value = "  +7.42%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "25.31"
$ws.Range("E46").Value = "  +12.30%  "
$ws.Range("B47").Value = "Filecoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.48"
$ws.Range("E47").Value = "  +4.67%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "41.70"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("E49").Value = "  +4.24%  "
$ws.Range("E50").Value = "  +6.96%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.378.16"
$ws.Range("E51").Value = "  +10.25%  "
